$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly price-history table: a new week (row 63, 2022-01-28) was added at the
# top of this block, pushing the existing rows 63-122 down by one to 64-123.
# Columns A,B,C,E,F,G,H,I,R are identical for every row in this market/product
# block, so only D (Fecha) and J,K,L,M,N,O,P,Q (the per-week stats) move.
$rows = @(
    @{ Row=63; D=44589; J=65; K=12000; L=12000; M=12000; N="`$/docena de atados"; O="Región de La Araucanía"; P=4000; Q=3 },
    @{ Row=64; D=44550; J=65; K=8000; L=8000; M=8000; N="`$/docena de atados"; O="Región de La Araucanía"; P=2667; Q=3 },
    @{ Row=65; D=44356; J=25; K=7000; L=8000; M=7600; N="`$/docena de atados"; O="Región de La Araucanía"; P=2533; Q=3 },
    @{ Row=66; D=44487; J=55; K=10000; L=10000; M=10000; N="`$/docena de atados"; O="Región de La Araucanía"; P=3333; Q=3 },
    @{ Row=67; D=44410; J=30; K=10000; L=10000; M=10000; N="`$/docena de atados"; O="Región de La Araucanía"; P=3333; Q=3 },
    @{ Row=68; D=44327; J=30; K=8000; L=8000; M=8000; N="`$/docena de atados"; O="Región de La Araucanía"; P=2667; Q=3 },
    @{ Row=69; D=44455; J=10; K=9000; L=9000; M=9000; N="`$/docena de atados"; O="Región de La Araucanía"; P=3000; Q=3 },
    @{ Row=70; D=44582; J=30; K=14000; L=14000; M=14000; N="`$/docena de atados"; O="Región de La Araucanía"; P=4667; Q=3 },
    @{ Row=71; D=44159; J=40; K=6000; L=6000; M=6000; N="`$/docena de atados"; O="Región de La Araucanía"; P=2000; Q=3 },
    @{ Row=72; D=44466; J=40; K=10000; L=10000; M=10000; N="`$/docena de atados"; O="Región de La Araucanía"; P=3333; Q=3 },
    @{ Row=73; D=44462; J=40; K=12000; L=12000; M=12000; N="`$/docena de atados"; O="Región de La Araucanía"; P=4000; Q=3 },
    @{ Row=74; D=44580; J=20; K=12000; L=12000; M=12000; N="`$/docena de atados"; O="Región de La Araucanía"; P=4000; Q=3 },
    @{ Row=75; D=44364; J=45; K=8000; L=8000; M=8000; N="`$/docena de atados"; O="Región de La Araucanía"; P=2667; Q=3 },
    @{ Row=76; D=44463; J=20; K=11000; L=12000; M=11500; N="`$/docena de atados"; O="Región de La Araucanía"; P=3833; Q=3 },
    @{ Row=77; D=44484; J=20; K=10000; L=10000; M=10000; N="`$/docena de atados"; O="Región de La Araucanía"; P=3333; Q=3 },
    @{ Row=78; D=44516; J=50; K=8000; L=8000; M=8000; N="`$/docena de atados"; O="Región de La Araucanía"; P=2667; Q=3 },
    @{ Row=79; D=44452; J=30; K=9000; L=9000; M=9000; N="`$/docena de atados"; O="Región de La Araucanía"; P=3000; Q=3 },
    @{ Row=80; D=44431; J=65; K=12000; L=12000; M=12000; N="`$/docena de atados"; O="Región de La Araucanía"; P=4000; Q=3 },
    @{ Row=81; D=44217; J=150; K=700; L=700; M=700; N="`$/docena de atados"; O="Región de La Araucanía"; P=233; Q=3 },
    @{ Row=82; D=44512; J=20; K=8000; L=8000; M=8000; N="`$/docena de atados"; O="Región de La Araucanía"; P=2667; Q=3 },
    @{ Row=83; D=44511; J=50; K=8000; L=8000; M=8000; N="`$/docena de atados"; O="Región de La Araucanía"; P=2667; Q=3 },
    @{ Row=84; D=44403; J=55; K=12000; L=12000; M=12000; N="`$/docena de atados"; O="Región de La Araucanía"; P=4000; Q=3 },
    @{ Row=85; D=44160; J=30; K=6000; L=6000; M=6000; N="`$/docena de atados"; O="Región de La Araucanía"; P=2000; Q=3 },
    @{ Row=86; D=44203; J=40; K=7000; L=7000; M=7000; N="`$/docena de atados"; O="Región de La Araucanía"; P=2333; Q=3 },
    @{ Row=87; D=44519; J=55; K=8000; L=8000; M=8000; N="`$/docena de atados"; O="Provincia de Chacabuco"; P=2667; Q=3 },
    @{ Row=88; D=44161; J=50; K=6000; L=6000; M=6000; N="`$/docena de atados"; O="Región de La Araucanía"; P=2000; Q=3 },
    @{ Row=89; D=44581; J=70; K=13000; L=14000; M=13571; N="`$/docena de atados"; O="Región de La Araucanía"; P=4524; Q=3 },
    @{ Row=90; D=44504; J=95; K=8000; L=8000; M=8000; N="`$/docena de atados"; O="Región de La Araucanía"; P=2667; Q=3 },
    @{ Row=91; D=44567; J=90; K=7000; L=9000; M=8111; N="`$/docena de atados"; O="Región de La Araucanía"; P=2704; Q=3 },
    @{ Row=92; D=44280; J=110; K=10000; L=10000; M=10000; N="`$/docena de atados"; O="Región de La Araucanía"; P=3333; Q=3 },
    @{ Row=93; D=44532; J=85; K=8000; L=8000; M=8000; N="`$/docena de atados"; O="Región de La Araucanía"; P=2667; Q=3 },
    @{ Row=94; D=44588; J=75; K=13000; L=14000; M=13533; N="`$/docena de atados"; O="Región de La Araucanía"; P=4511; Q=3 },
    @{ Row=95; D=44442; J=10; K=9000; L=9000; M=9000; N="`$/docena de atados"; O="Región de La Araucanía"; P=3000; Q=3 },
    @{ Row=96; D=44344; J=40; K=8000; L=8000; M=8000; N="`$/docena de atados"; O="Región de La Araucanía"; P=2667; Q=3 },
    @{ Row=97; D=44334; J=25; K=7000; L=7000; M=7000; N="`$/docena de atados"; O="Región de La Araucanía"; P=2333; Q=3 },
    @{ Row=98; D=44266; J=65; K=10000; L=10000; M=10000; N="`$/docena de atados"; O="Región de La Araucanía"; P=3333; Q=3 },
    @{ Row=99; D=44277; J=65; K=10000; L=10000; M=10000; N="`$/docena de atados"; O="Región de La Araucanía"; P=3333; Q=3 },
    @{ Row=100; D=44454; J=30; K=9000; L=9000; M=9000; N="`$/docena de atados"; O="Región de La Araucanía"; P=3000; Q=3 },
    @{ Row=101; D=44371; J=30; K=8000; L=8000; M=8000; N="`$/docena de atados"; O="Región de La Araucanía"; P=2667; Q=3 },
    @{ Row=102; D=44259; J=40; K=11000; L=11000; M=11000; N="`$/docena de atados"; O="Región de La Araucanía"; P=3667; Q=3 },
    @{ Row=103; D=44330; J=20; K=8000; L=8000; M=8000; N="`$/docena de atados"; O="Región de La Araucanía"; P=2667; Q=3 },
    @{ Row=104; D=44494; J=20; K=8000; L=8000; M=8000; N="`$/docena de atados"; O="Región de La Araucanía"; P=2667; Q=3 },
    @{ Row=105; D=44526; J=40; K=8000; L=8000; M=8000; N="`$/docena de atados"; O="Región de La Araucanía"; P=2667; Q=3 },
    @{ Row=106; D=44533; J=110; K=8000; L=8000; M=8000; N="`$/docena de atados"; O="Región de La Araucanía"; P=2667; Q=3 },
    @{ Row=107; D=44354; J=80; K=7000; L=8000; M=7500; N="`$/docena de atados"; O="Región de La Araucanía"; P=2500; Q=3 },
    @{ Row=108; D=44221; J=90; K=7000; L=8000; M=7611; N="`$/docena de atados"; O="Región de La Araucanía"; P=2537; Q=3 },
    @{ Row=109; D=44523; J=30; K=8000; L=8000; M=8000; N="`$/docena de atados"; O="Región de La Araucanía"; P=2667; Q=3 },
    @{ Row=110; D=44399; J=20; K=12000; L=12000; M=12000; N="`$/docena de atados"; O="Región de La Araucanía"; P=4000; Q=3 },
    @{ Row=111; D=44382; J=50; K=8000; L=9000; M=8600; N="`$/docena de atados"; O="Región de La Araucanía"; P=2867; Q=3 },
    @{ Row=112; D=44441; J=50; K=9000; L=10000; M=9600; N="`$/docena de atados"; O="Región de La Araucanía"; P=3200; Q=3 },
    @{ Row=113; D=44372; J=30; K=8000; L=8000; M=8000; N="`$/docena de atados"; O="Región de La Araucanía"; P=2667; Q=3 },
    @{ Row=114; D=44335; J=35; K=7000; L=7000; M=7000; N="`$/docena de atados"; O="Región de La Araucanía"; P=2333; Q=3 },
    @{ Row=115; D=44438; J=30; K=9000; L=9000; M=9000; N="`$/docena de atados"; O="Región de La Araucanía"; P=3000; Q=3 },
    @{ Row=116; D=44453; J=20; K=9000; L=9000; M=9000; N="`$/docena de atados"; O="Región de La Araucanía"; P=3000; Q=3 },
    @{ Row=117; D=44385; J=40; K=9000; L=10000; M=9500; N="`$/docena de atados"; O="Región de La Araucanía"; P=3167; Q=3 },
    @{ Row=118; D=44162; J=50; K=6000; L=6000; M=6000; N="`$/docena de atados"; O="Región de La Araucanía"; P=2000; Q=3 },
    @{ Row=119; D=44529; J=65; K=10000; L=10000; M=10000; N="`$/cuna 10 kilos"; O="Región Metropolitana"; P=1000; Q=10 },
    @{ Row=120; D=44299; J=70; K=9000; L=9000; M=9000; N="`$/docena de atados"; O="Región de La Araucanía"; P=3000; Q=3 },
    @{ Row=121; D=44428; J=10; K=8000; L=8000; M=8000; N="`$/docena de atados"; O="Región de La Araucanía"; P=2667; Q=3 },
    @{ Row=122; D=44302; J=40; K=8000; L=8000; M=8000; N="`$/docena de atados"; O="Región de La Araucanía"; P=2667; Q=3 },
    @{ Row=123; D=44209; J=50; K=7000; L=7000; M=7000; N="`$/docena de atados"; O="Región de La Araucanía"; P=2333; Q=3 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 10).Value = $r.J
    $ws.Cells.Item($r.Row, 11).Value = $r.K
    $ws.Cells.Item($r.Row, 12).Value = $r.L
    $ws.Cells.Item($r.Row, 13).Value = $r.M
    $ws.Cells.Item($r.Row, 14).Value = $r.N
    $ws.Cells.Item($r.Row, 15).Value = $r.O
    $ws.Cells.Item($r.Row, 16).Value = $r.P
    $ws.Cells.Item($r.Row, 17).Value = $r.Q
}

# Row 123 is brand new; give its date cell (D123) the same date/time number
# format used by the rest of column D (style index 2 in styles.xml), then
# fill in the row-constant columns that are shared by the whole block.
$ws.Range("D123").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(123, 1).Value = 10
$ws.Cells.Item(123, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(123, 3).Value = "La Araucanía"
$ws.Cells.Item(123, 5).Value = 9
$ws.Cells.Item(123, 6).Value = 100112012
$ws.Cells.Item(123, 7).Value = "Espinaca"
$ws.Cells.Item(123, 8).Value = "Sin especificar"
$ws.Cells.Item(123, 9).Value = "Primera"
$ws.Cells.Item(123, 18).Value = "Hortaliza"
